$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new "venv" config entry -----------------------------------
# Column B: path to the venv activation script, with "booking-venv"
# rendered in bold navy to highlight the venv name within the path.
$full = "D:\bookings\booking-venv\Scripts\activate.ps1"
$ws.Range("B7").Value = $full

$boldStart = $full.IndexOf("booking-venv") + 1
$boldLen = "booking-venv".Length
$boldChars = $ws.Range("B7").Characters($boldStart, $boldLen)
$boldChars.Font.Bold = $true
$boldChars.Font.Color = 0x701919
$boldChars.Font.Name = "Calibri"
$boldChars.Font.Size = 11

$tailStart = $boldStart + $boldLen
$tailLen = $full.Length - $tailStart + 1
$tailChars = $ws.Range("B7").Characters($tailStart, $tailLen)
$tailChars.Font.Color = 0x000000
$tailChars.Font.Name = "Calibri"
$tailChars.Font.Size = 11

# Column C: description
$ws.Range("C7").Value = "Path to .pst activate script for venv"

# Column A: name (set last to mirror the shared-string insertion order)
$ws.Range("A7").Value = "venv"

# Match row height of the other data rows
$ws.Rows.Item(7).RowHeight = 15

# Leave selection where Excel would after the edit
$ws.Range("C8").Select()

Write-Output "Row 7 (venv) added"
